# Update reports - 2026-01-30 09:55
# Adds the new Welsh Affairs Committee government-response report to the
# "Reports" table, and records the scan that discovered it in the "Scans"
# table.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "Reports" sheet/table: append the new publication as row 33
# ---------------------------------------------------------------------
$wsReports = $wb.Worksheets.Item("Reports")
$reportsTable = $wsReports.ListObjects.Item("Reports")
$reportsTable.ListRows.Add() | Out-Null

# Publication ID and Publication Date look like plain numbers/dates, so a
# leading apostrophe is used to force them to stay as text (matching the
# rest of the sheet), then the style is reset back to Normal so no extra
# cell formatting is left behind.
$wsReports.Range("A33").Value = "'51312"
$wsReports.Range("A33").Style = "Normal"
$wsReports.Range("B33").Value = "HC 785"
$wsReports.Range("C33").Value = "2024-26"
$wsReports.Range("D33").Value = "Welsh Affairs Committee"
$wsReports.Range("E33").Value = "Commons"
$wsReports.Range("F33").Value = "Farming in Wales in 2025: Challenges and Opportunities: Government Response"
$wsReports.Range("G33").Value = "1st Special Report"
$wsReports.Range("H33").Value = "'2026-01-30"
$wsReports.Range("H33").Style = "Normal"
$wsReports.Range("I33").Value = "09:00:00"
# J33 (Late by min) intentionally left blank
$wsReports.Range("K33").Value = "0:55:18"

# ---------------------------------------------------------------------
# 2) "Scans" sheet/table: append the scan record as row 17
# ---------------------------------------------------------------------
$wsScans = $wb.Worksheets.Item("Scans")
$scansTable = $wsScans.ListObjects.Item("Scans")
$scansTable.ListRows.Add() | Out-Null

$wsScans.Range("A17").Value = "'2026-01-30"
$wsScans.Range("A17").Style = "Normal"
$wsScans.Range("B17").Value = "09:55:18"
$wsScans.Range("C17").Value = "'51312"
$wsScans.Range("C17").Style = "Normal"
